# Generate Report for Handback
# Update status text (shared across Overview + per-language sheets) and
# populate the handback columns (Latest Target File / Latest Handback File /
# Latest Handback DateTime) now that the localized files are in sync.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f39cd4c9debd6a6bea08b64d842b5393f3930a15/e2e/a.md"

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace($oldStatus, $newStatus) | Out-Null
}

function Set-HandbackColumns($ws, $xlfName, $handbackDateTime) {
    foreach ($row in 2, 3) {
        $targetCell = $ws.Range("I$row")
        $targetCell.Value = "a.md"
        $ws.Hyperlinks.Add($targetCell, $baseUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
        $targetCell.Font.Underline = $true
        $targetCell.Font.Color = 15570276  # Cornflower blue (FF6495ED), matches the Source File Name hyperlink style

        $ws.Range("J$row").Value = $xlfName
        $ws.Range("K$row").Value = $handbackDateTime
    }
}

# zh-cn sheet
$zhcn = $wb.Worksheets.Item("zh-cn")
Set-HandbackColumns $zhcn "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-08-16 02:33:36"

# de-de sheet
$dede = $wb.Worksheets.Item("de-de")
Set-HandbackColumns $dede "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-08-16 02:33:42"
